$d = $word.ActiveDocument

# Locate the run-spanning text "<id>p167r_1</id>" which today is split across
# three separate runs (the first/last carrying Courier-New/7f6000 formatting,
# the middle one plain). We want to collapse it into a single run that keeps
# the first run's formatting and contains the whole string as its text.
$full = $d.Content
$found = $full.Find.Execute("<id>p167r_1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $full.Start
    $end = $full.End

    # Keep the leading "<id>" (first run, 4 characters) exactly as-is so its
    # run formatting (Courier New / color 7f6000 / sz 18) is preserved, and
    # remove everything after it in this match ("p167r_1</id>", currently the
    # 2nd + 3rd runs).
    $idLen = 4
    $tailStart = $start + $idLen

    $tailRng = $d.Range($tailStart, $end)
    $tailRng.Delete()

    # Re-insert the removed text right after the kept "<id>" run; inserting at
    # that boundary with nothing else selected makes Word continue the
    # preceding run's formatting, merging everything into one run.
    $insertRng = $d.Range($tailStart, $tailStart)
    $insertRng.InsertAfter("p167r_1</id>")
}
